# Add an "Address" column before the existing "District" column (F).
# This shifts the existing column F ("District") to column G, and
# populates the new column F with address information for each teacher.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at F; this pushes the existing F column (District) to G.
$ws.Columns.Item(6).Insert()

# Header
$ws.Range("F2").Value = "Address"

# Address values extracted for each row (blank where the source row had no
# distinguishable address component of its own).
$addresses = @{
    3  = "G H S RevatagaonIndi"
    4  = "Govt. High School Kambagi"
    5  = "Karnataka High SchoolBurnapur Road"
    6  = "S S H S GajevaniIndi"
    7  = "S G R J D H S Nimbal RSIndi"
    8  = "S S H S TambaIndi"
    9  = "S M High School MulawadB Bagewadi"
    10 = "Govt. P U College for Boys(High School section)"
    11 = "T S S H S Muddebihal"
    12 = "M D S P H S NarasalagB Bagewadi"
    13 = "GHS GaniB. Bagewadi"
    14 = "N E H S NidagundiB Bagewadi"
    15 = "G H P S Nagaral"
    16 = "Govt. P U College for Girls"
    18 = "M D R S LachyanIndi"
    19 = "Shri Satya Saibaba High School Ghonasags"
    20 = "High School"
    21 = "D N Darbar Govt. High School"
    22 = "Shantala H S Muddebihal"
    23 = "Govt. S B P U College( High School Section)B Bagewadi"
    24 = "G H S RakkasagiMuddebihal"
    25 = "G H S TambaIndi"
    27 = "Shri Hanuman High School"
    28 = "J J H S WadawadagiB Bagewadi"
    29 = "Govt. High School devaragennur"
    30 = "Govt. High School UmaraniIndi"
    31 = "G H P S Makhanapur LT-1"
    32 = "Ratnapur,Tajapur(H)"
}

foreach ($row in $addresses.Keys) {
    $ws.Cells.Item($row, 6).Value = $addresses[$row]
}

# Rows 17 and 26 have no address text of their own, so column F must remain blank there.
$ws.Cells.Item(17, 6).Value = ""
$ws.Cells.Item(26, 6).Value = ""
